# Apply crypto price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-valued cells (names, links, percent strings, multi-dot prices)
# and price strings that are not parsed as pure numbers by Excel.
$textUpdates = @{
    "D2" = "94.467.29"
    "E2" = "  +1.74%  "
    "D3" = "3.077.23"
    "E3" = "  -0.95%  "
    "E4" = "  -0.16%  "
    "E5" = "  -2.33%  "
    "E6" = "  -0.29%  "
    "E7" = "  -1.25%  "
    "E8" = "  -3.61%  "
    "E9" = "  +0.03%  "
    "E10" = "  +7.27%  "
    "D11" = "3.071.45"
    "E11" = "  -1.12%  "
    "E12" = "  -2.56%  "
    "D13" = "94.080.34"
    "E13" = "  +1.17%  "
    "E14" = "  -3.77%  "
    "E15" = "  -1.45%  "
    "E16" = "  -2.29%  "
    "D17" = "3.645.05"
    "E17" = "  -1.58%  "
    "D18" = "3.040.87"
    "E18" = "  -2.37%  "
    "E19" = "  -5.74%  "
    "E20" = "  -2.31%  "
    "E21" = "  -1.60%  "
    "E22" = "  -1.29%  "
    "E23" = "  -4.87%  "
    "E24" = "  -5.60%  "
    "E25" = "  +5.95%  "
    "E26" = "  -4.32%  "
    "E27" = "  -2.14%  "
    "E28" = "  +2.26%  "
    "D29" = "3.226.69"
    "E30" = "  +0.09%  "
    "E31" = "  +6.69%  "
    "E32" = "  +5.07%  "
    "E33" = "  -8.57%  "
    "E34" = "  -0.81%  "
    "E35" = "  -4.83%  "
    "E36" = "  -2.63%  "
    "B37" = "EthereumClassic"
    "C37" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "E37" = "  -2.16%  "
    "B38" = "Binance-PegBSC-USD"
    "C38" = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
    "E38" = "  +0.80%  "
    "E39" = "  -1.08%  "
    "E40" = "  +3.51%  "
    "E41" = "  +1.55%  "
    "E42" = "  -4.28%  "
    "E43" = "  -3.64%  "
    "E44" = "  -2.12%  "
    "E45" = "  -0.02%  "
    "E46" = "  -7.79%  "
    "E47" = "  -1.01%  "
    "E48" = "  -2.02%  "
    "E49" = "  -4.03%  "
    "E50" = "  -1.13%  "
    "E51" = "  -0.13%  "
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Price cells whose new text happens to look like a plain number (e.g. "235.01").
# Excel auto-converts Range.Value assignments of numeric-looking strings into
# real floating point numbers (losing the original text formatting / precision,
# e.g. "14.40" -> 14.4, or introducing binary float noise). The source workbook
# stores these as literal text, so force text entry by temporarily switching the
# cell to a Text number format, assign the value, then restore the original style
# so no visible formatting change is left behind.
$numericTextUpdates = @{
    "D5" = "235.01"
    "D6" = "609.22"
    "D8" = "0.378"
    "D10" = "0.805"
    "D12" = "0.196"
    "D14" = "0.0000240"
    "D15" = "33.73"
    "D16" = "5.32"
    "D19" = "3.55"
    "D20" = "14.40"
    "D21" = "5.66"
    "D22" = "439.77"
    "D23" = "8.81"
    "D24" = "0.0000189"
    "D26" = "5.50"
    "D27" = "84.69"
    "D28" = "11.90"
    "D31" = "0.247"
    "D34" = "9.04"
    "D35" = "7.63"
    "D37" = "25.41"
    "D38" = "0.886"
    "D41" = "0.436"
    "D42" = "468.57"
    "D43" = "3.70"
    "D44" = "1.26"
    "D46" = "3.09"
    "D47" = "161.60"
    "D48" = "0.669"
    "D50" = "43.54"
    "D51" = "0.997"
}
foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = $origStyle
}
